$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheets index 1 / sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2232
$ws1.Range("F5").Value = 13439
$ws1.Range("F8").Value = 526
$ws1.Range("F9").Value = 491
$ws1.Range("F11").Value = 1009
$ws1.Range("F12").Value = 13832
$ws1.Range("F13").Value = 14529
$ws1.Range("F22").Value = 48
$ws1.Range("F23").Value = 6
$ws1.Range("F24").Value = 1113
$ws1.Range("F27").Value = 5567
$ws1.Range("F29").Value = 1039
$ws1.Range("F31").Value = 36
$ws1.Range("F32").Value = 24
$ws1.Range("F33").Value = 154

# Sheet "演出" (Worksheets index 2 / sheet2.xml)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 2

# Sheet "全部类型" (Worksheets index 4 / sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2232
$ws4.Range("F5").Value = 13439
$ws4.Range("F7").Value = 2
$ws4.Range("F9").Value = 526
$ws4.Range("F10").Value = 491
$ws4.Range("F12").Value = 1009
$ws4.Range("F13").Value = 13832
$ws4.Range("F14").Value = 14529
$ws4.Range("F23").Value = 48
$ws4.Range("F24").Value = 6
$ws4.Range("F25").Value = 1113
$ws4.Range("F28").Value = 5567
$ws4.Range("F30").Value = 1039
$ws4.Range("F32").Value = 36
$ws4.Range("F33").Value = 24
$ws4.Range("F34").Value = 154
